$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Value  = 0.1067
$ws.Range("I6").Value  = -0.0452
$ws.Range("J6").Value  = 0.295
$ws.Range("K6").Value  = 0.2573
$ws.Range("L6").Value  = 0.3088
$ws.Range("M6").Value  = 0.2796
$ws.Range("N6").Value  = 0.1561
$ws.Range("O6").Value  = -0.0145
$ws.Range("P6").Value  = -0.0203
$ws.Range("Q6").Value  = -0.0138
$ws.Range("R6").Value  = -0.0083
$ws.Range("S6").Value  = 0.0126
$ws.Range("T6").Value  = 0.0046
$ws.Range("U6").Value  = -0.005
$ws.Range("V6").Value  = -0.0105
$ws.Range("W6").Value  = -1.7941

$ws.Range("H16").Value = -2.2874
$ws.Range("I16").Value = -0.5607
$ws.Range("J16").Value = 0.0638
$ws.Range("K16").Value = -0.1884
$ws.Range("L16").Value = 0.5287
$ws.Range("M16").Value = 0.2312
$ws.Range("N16").Value = -0.2545
$ws.Range("O16").Value = -0.6538
$ws.Range("P16").Value = -0.2369
$ws.Range("Q16").Value = -0.0086
$ws.Range("R16").Value = -0.5226
$ws.Range("S16").Value = -0.6951
$ws.Range("T16").Value = -0.4849
$ws.Range("U16").Value = -0.4247
$ws.Range("V16").Value = -0.0558
$ws.Range("W16").Value = -73.1408

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("O34").Value = 0.1058
$ws.Range("P34").Value = 0.2178
$ws.Range("Q34").Value = 0.2097
$ws.Range("R34").Value = 0.2755
$ws.Range("S34").Value = 0.2321
$ws.Range("T34").Value = 0.1752
$ws.Range("U34").Value = 0.0305
$ws.Range("V34").Value = 0.1036
$ws.Range("W34").Value = -1.3669

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = 0
$ws.Range("N44").Value = 0
$ws.Range("O44").Value = 0.02
$ws.Range("P44").Value = 0.1261
$ws.Range("Q44").Value = 0.1105
$ws.Range("R44").Value = 0.1636
$ws.Range("S44").Value = 0.1244
$ws.Range("T44").Value = 0.0793
$ws.Range("U44").Value = -0.0543
$ws.Range("V44").Value = 0.0352
$ws.Range("W44").Value = -2.4937
